# Update CF for H (Hydrogen) and drop a stale duplicate "Methane / urban" row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Hydrogen (rows 120-124): characterization factor 4.3 -> 5
$ws.Range("C120:C124").Value = 5

# 2) Row 125 ("Methane" / "air::urban air close to ground" / 4.3) was a stale
#    duplicate of the row that followed it (same name/category, value 29.7).
#    Delete it so the table collapses back to one row per name/category and
#    everything below shifts up by one.
$ws.Rows(125).Delete()

# Leave the selection where the editor left it.
$ws.Range("B112").Select()
